{"js": "// Apply the \"roach 2 manufacture updated\" edit.\n//\n// Summary of the change:\n//  1. \"The manufacturing pack consists of x folders\" ->\n//     \"The manufacturing pack consists of 4 folders:\"\n//  2. Four new list paragraphs (+ one blank paragraph) are inserted right\n//     after that paragraph, describing the folders that ship in the pack.\n//  3. \"  X BACK PANEL\" -> \"2 X BACK PANEL, blanks by default\"\n//  4. \"...brushed anodised finished\" -> \"...brushed anodised finish.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1: locate the \"consists of x folders\" paragraph ---------------------\nconst introIndex = paragraphs.items.findIndex(p =>\n  p.text.indexOf(\"manufacturing pack consists of\") !== -1\n);\nif (introIndex === -1) {\n  throw new Error(\"Could not find the 'manufacturing pack consists of' paragraph\");\n}\nconst introPara = paragraphs.items[introIndex];\n\n// Update \"x folders\" -> \"4 folders:\" (keeps the leading \"The \" text intact).\nconst introSearch = introPara.search(\"x folders\", { matchCase: true });\nintroSearch.load(\"items\");\nawait context.sync();\nif (introSearch.items.length === 0) {\n  throw new Error(\"Could not find 'x folders' text to replace\");\n}\nintroSearch.items[0].insertText(\"4 folders:\", \"Replace\");\nawait context.sync();\n\n// --- 2: insert the folder-description paragraphs + trailing blank line ---\n// Inserting a single block of text with embedded \"\\n\"s (rather than calling\n// insertParagraph() repeatedly) makes the final blank paragraph come out as\n// a true empty <w:p/> instead of a paragraph holding an empty run.\nconst newLines = [\n  \"dcoumentation \\u2013 Drawings and Information\",\n  \"dxf \\u2013 Autocad 2d files\",\n  \"silkscreens \\u2013 Silkscreens for the front panel, back panels and back\",\n  \"stp_files \\u2013 3D Step Files\",\n  \"\"\n];\nconst introEnd = introPara.getRange(\"End\");\nintroEnd.insertText(\"\\n\" + newLines.join(\"\\n\"), \"End\");\nawait context.sync();\n\n// --- 3: \"  X BACK PANEL\" -> \"2 X BACK PANEL, blanks by default\" ----------\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst backPanelIndex = paragraphs.items.findIndex(p =>\n  p.text.indexOf(\"X BACK PANEL\") !== -1\n);\nif (backPanelIndex === -1) {\n  throw new Error(\"Could not find the 'X BACK PANEL' paragraph\");\n}\nconst backPanelPara = paragraphs.items[backPanelIndex];\nbackPanelPara.insertText(\"2 X BACK PANEL, blanks by default\", \"Replace\");\nawait context.sync();\n\n// --- 4: \"...anodised finished\" -> \"...anodised finish.\" ------------------\nconst finishSearch = body.search(\"ed finished\", { matchCase: true });\nfinishSearch.load(\"items\");\nawait context.sync();\nif (finishSearch.items.length === 0) {\n  throw new Error(\"Could not find 'ed finished' text to replace\");\n}\nfinishSearch.items[0].insertText(\"ed finish.\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1: \"x folders\" -> \"4 folders:\" ---------------------------------------\n$find = $d.Content.Find\n$find.Text = \"x folders\"\n$find.Replacement.Text = \"4 folders:\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- 2: insert the folder-description paragraphs + trailing blank line ---\n$introPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*manufacturing pack consists of*\") {\n        $introPara = $p\n        break\n    }\n}\n$anchor = $introPara.Range\n$anchor.InsertParagraphAfter()\n$newParaIndex = $introPara.Index + 1\n$lines = \"dcoumentation \" + [char]0x2013 + \" Drawings and Information`r\" + `\n         \"dxf \" + [char]0x2013 + \" Autocad 2d files`r\" + `\n         \"silkscreens \" + [char]0x2013 + \" Silkscreens for the front panel, back panels and back`r\" + `\n         \"stp_files \" + [char]0x2013 + \" 3D Step Files`r\"\n$d.Paragraphs($newParaIndex).Range.Text = $lines\n\n# --- 3: \"  X BACK PANEL\" -> \"2 X BACK PANEL, blanks by default\" ----------\n$backPanelPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*X BACK PANEL*\") {\n        $backPanelPara = $p\n        break\n    }\n}\n$backPanelPara.Range.Text = \"2 X BACK PANEL, blanks by default\"\n\n# --- 4: \"...anodised finished\" -> \"...anodised finish.\" ------------------\n$find2 = $d.Content.Find\n$find2.Text = \"ed finished\"\n$find2.Replacement.Text = \"ed finish.\"\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
